$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" forces Excel to treat the assignment as literal text so COM
# does not auto-convert numeric-looking strings (e.g. "1.00", "0.0510",
# "34.428.98") into real numbers -- losing the original text formatting --
# while leaving each cell's number format / style untouched. The leading
# quote itself is a formatting marker and is not stored in the cell value.
$ws.Range("D2").Value = "'34.428.98"
$ws.Range("E2").Value = "'  +0.45%  "
$ws.Range("D3").Value = "'1.799.27"
$ws.Range("E3").Value = "'  +0.49%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.23%  "
$ws.Range("D5").Value = "'224.54"
$ws.Range("E5").Value = "'  +0.01%  "
$ws.Range("E6").Value = "'  +1.41%  "
$ws.Range("E7").Value = "'  +0.32%  "
$ws.Range("D8").Value = "'41.19"
$ws.Range("E8").Value = "'  +14.59%  "
$ws.Range("E9").Value = "'  +0.13%  "
$ws.Range("E10").Value = "'  -0.91%  "
$ws.Range("D11").Value = "'0.0997"
$ws.Range("E11").Value = "'  +3.77%  "
$ws.Range("D12").Value = "'2.059.71"
$ws.Range("E12").Value = "'  +0.39%  "
$ws.Range("D13").Value = "'1.798.73"
$ws.Range("E13").Value = "'  +0.32%  "
$ws.Range("E14").Value = "'  -2.48%  "
$ws.Range("D15").Value = "'34.439.31"
$ws.Range("E15").Value = "'  +0.43%  "
$ws.Range("E16").Value = "'  -0.44%  "
$ws.Range("D17").Value = "'4.37"
$ws.Range("E17").Value = "'  +0.77%  "
$ws.Range("D18").Value = "'67.21"
$ws.Range("E18").Value = "'  -1.84%  "
$ws.Range("D19").Value = "'239.36"
$ws.Range("E19").Value = "'  -0.01%  "
$ws.Range("E20").Value = "'  -0.18%  "
$ws.Range("D21").Value = "'11.06"
$ws.Range("E21").Value = "'  -0.97%  "
$ws.Range("E22").Value = "'  +0.30%  "
$ws.Range("D23").Value = "'4.09"
$ws.Range("E23").Value = "'  +1.18%  "
$ws.Range("E24").Value = "'  -0.94%  "
$ws.Range("D25").Value = "'171.33"
$ws.Range("E25").Value = "'  +0.62%  "
$ws.Range("D26").Value = "'7.61"
$ws.Range("E26").Value = "'  -3.60%  "
$ws.Range("E27").Value = "'  +1.03%  "
$ws.Range("E28").Value = "'  +0.56%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "'  +0.19%  "
$ws.Range("D30").Value = "'1.23"
$ws.Range("E30").Value = "'  +0.24%  "
$ws.Range("E31").Value = "'  +0.61%  "
$ws.Range("D32").Value = "'3.84"
$ws.Range("E32").Value = "'  -0.41%  "
$ws.Range("D33").Value = "'0.0510"
$ws.Range("E33").Value = "'  +0.17%  "
$ws.Range("E34").Value = "'  +0.48%  "
$ws.Range("D35").Value = "'1.315.53"
$ws.Range("E35").Value = "'  -2.93%  "
$ws.Range("D36").Value = "'0.641"
$ws.Range("E36").Value = "'  +0.29%  "
$ws.Range("E37").Value = "'  +0.91%  "
$ws.Range("D38").Value = "'85.03"
$ws.Range("E38").Value = "'  +6.16%  "
$ws.Range("B39").Value = "'VeChain"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0187"
$ws.Range("E39").Value = "'  +1.89%  "
$ws.Range("B40").Value = "'RenderToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.33"
$ws.Range("E40").Value = "'  +0.48%  "
$ws.Range("B41").Value = "'InjectiveProtocol"
$ws.Range("C41").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'14.65"
$ws.Range("E41").Value = "'  +12.26%  "
$ws.Range("B42").Value = "'WEMIXToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.24"
$ws.Range("E42").Value = "'  +6.00%  "
$ws.Range("E43").Value = "'  +0.91%  "
$ws.Range("E44").Value = "'  +0.14%  "
$ws.Range("D45").Value = "'0.934"
$ws.Range("E45").Value = "'  +1.05%  "
$ws.Range("D46").Value = "'0.0519"
$ws.Range("E46").Value = "'  +4.52%  "
$ws.Range("D47").Value = "'1.960.24"
$ws.Range("E47").Value = "'  +0.36%  "
$ws.Range("D48").Value = "'5.83"
$ws.Range("E48").Value = "'  +1.47%  "
$ws.Range("E49").Value = "'  +0.31%  "
$ws.Range("D50").Value = "'100.59"
$ws.Range("E50").Value = "'  -0.63%  "
$ws.Range("E51").Value = "'  +1.61%  "
